# Refresh the coinranking.com "cryptos" snapshot: Price (D) and
# Volume(1h) (E) columns for every coin row (2-51) are overwritten
# with the latest scrape, as published by the scheduled GitHub
# Actions job ("Updated cryptos list ... with GitHub Actions").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry is the target cell plus its new display text. ForceText
# is set for Price values that are plain digits/dots (e.g. "1.002",
# "223.41") - Excel would otherwise auto-convert those into numbers
# on assignment (dropping significant trailing zeros such as in
# "7.200" -> 7.2, or mis-parsing "1.040.29" as a date/number), so we
# pre-format the cell as Text to keep it an exact string, matching
# how the other Price cells (already non-numeric, e.g. "27.268.30")
# are stored.
$updates = @(
    @{ Cell = 'D2'; Value = '27.268.30'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -0.81%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.701.83'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -1.30%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '1.002'; ForceText = $true },
    @{ Cell = 'D5'; Value = '223.41'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -1.12%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '0.5305'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -1.38%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '1.002'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -0.11%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.2659'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -0.87%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.06584'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -0.27%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '20.71'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  -4.55%  '; ForceText = $false },
    @{ Cell = 'E11'; Value = '  -1.45%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '4.490'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  -3.17%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '1.717.98'; ForceText = $false },
    @{ Cell = 'E13'; Value = '  -0.39%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '1.936.62'; ForceText = $false },
    @{ Cell = 'E14'; Value = '  -1.25%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '0.5782'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -1.62%  '; ForceText = $false },
    @{ Cell = 'E16'; Value = '  -1.68%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '67.41'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -0.92%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '27.262.29'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  -0.89%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '215.05'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -3.23%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '1.002'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  -0.11%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '4.606'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -2.74%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '10.33'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -3.28%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '5.965'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -2.23%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '1.003'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -0.15%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '143.90'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -2.72%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '1.702'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  +0.75%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '0.1197'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -2.74%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '7.198'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -2.89%  '; ForceText = $false },
    @{ Cell = 'E29'; Value = '  -3.50%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '0.05361'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -3.53%  '; ForceText = $false },
    @{ Cell = 'E31'; Value = '  -1.52%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '3.463'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -2.44%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '3.397'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -2.08%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '1.638'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -1.19%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  +1.81%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '2.414'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -1.33%  '; ForceText = $false },
    @{ Cell = 'E37'; Value = '  -1.57%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.5803'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -2.01%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '0.01628'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -0.93%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '5.773'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -1.52%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '1.002'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -0.12%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '1.039.44'; ForceText = $false },
    @{ Cell = 'E42'; Value = '  -1.52%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '0.8383'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  -2.09%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '100.92'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -0.77%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '1.844.22'; ForceText = $false },
    @{ Cell = 'E45'; Value = '  -1.22%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '0.0₈116'; ForceText = $false },
    @{ Cell = 'E46'; Value = '  +0.79%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '57.70'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -2.02%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '0.4515'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  +1.70%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '1.004'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  +0.43%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '8.017'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -2.41%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.05226'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -0.93%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
